# Updates cryptos list data in the active worksheet
# (price/volume refresh + two coin rows reordered, per upstream data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric must be forced back to
# Text format first, since the sheet stores all Price values as strings
# (e.g. thousands separated with '.' like '69.704.96').
$textCells = @("D5","D6","D13","D15","D17","D19","D21","D22","D23","D25","D26","D27","D29","D31","D32","D36","D37","D39","D42","D44","D47","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '69.704.96'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '2.511.39'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '575.89'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = '166.89'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '2.510.65'
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("E12").Value = '  +3.89%  '
$ws.Range("D13").Value = '4.93'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").Value = '2.969.87'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000179'
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '69.570.73'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '24.91'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '2.511.05'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '11.29'
$ws.Range("E19").Value = '  -1.30%  '
$ws.Range("E20").Value = '  -2.86%  '
$ws.Range("D21").Value = '350.49'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '3.92'
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").Value = '1.95'
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '70.23'
$ws.Range("E25").Value = '  +2.06%  '
$ws.Range("D26").Value = '3.96'
$ws.Range("E26").Value = '  -1.47%  '
$ws.Range("D27").Value = '8.85'
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").Value = '2.637.00'
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").Value = '0.0₃0893'
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("D32").Value = '462.40'
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("E33").Value = '  -4.94%  '
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '158.80'
$ws.Range("E36").Value = '  +2.11%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.117'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").Value = '18.50'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = '4.70'
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("D44").Value = '38.15'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E45").Value = '  -4.41%  '
$ws.Range("E46").Value = '  -8.00%  '
$ws.Range("D47").Value = '142.50'
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").Value = '0.521'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("E51").Value = '  -1.06%  '
